$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 104.23077
$ws.Cells.Item(5, 9).Value = 90.5
$ws.Cells.Item(5, 10).Value = 179.75
$ws.Cells.Item(5, 11).Value = 90.5
$ws.Cells.Item(5, 12).Value = 179.75
$ws.Cells.Item(5, 13).Value = 24.5
$ws.Cells.Item(5, 14).Value = -409.75
$ws.Cells.Item(12, 8).Value = 695.25
$ws.Cells.Item(12, 9).Value = 606.7143
$ws.Cells.Item(12, 11).Value = 606.7143
$ws.Cells.Item(12, 13).Value = -436.7143
$ws.Cells.Item(17, 8).Value = 1300.5834
$ws.Cells.Item(17, 10).Value = 1300.5834
$ws.Cells.Item(17, 12).Value = 3901.7502
$ws.Cells.Item(17, 14).Value = -4237.7502
$ws.Cells.Item(28, 8).Value = 1717.7778
$ws.Cells.Item(28, 9).Value = 1388.75
$ws.Cells.Item(28, 11).Value = 1388.75
$ws.Cells.Item(28, 13).Value = -903.75
$ws.Cells.Item(29, 8).Value = 571
$ws.Cells.Item(29, 9).Value = 571
$ws.Cells.Item(29, 11).Value = 1713
$ws.Cells.Item(29, 13).Value = -1432
$ws.Cells.Item(32, 8).Value = 0
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 14).Value = $null
$ws.Cells.Item(51, 8).Value = 3892.75
$ws.Cells.Item(51, 9).Value = 3690.3333
$ws.Cells.Item(51, 10).Value = 4500
$ws.Cells.Item(51, 11).Value = 3690.3333
$ws.Cells.Item(51, 12).Value = 4500
$ws.Cells.Item(51, 13).Value = -3206.3333
$ws.Cells.Item(51, 14).Value = -5468
$ws.Cells.Item(58, 8).Value = 488.3
$ws.Cells.Item(58, 9).Value = 172.875
$ws.Cells.Item(58, 10).Value = 1750
$ws.Cells.Item(58, 11).Value = 518.625
$ws.Cells.Item(58, 12).Value = 5250
$ws.Cells.Item(58, 13).Value = -368.625
$ws.Cells.Item(58, 14).Value = -5550
$ws.Cells.Item(62, 8).Value = 0
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 13).Value = $null
$ws.Cells.Item(64, 8).Value = 6670.8486
$ws.Cells.Item(64, 9).Value = 6504.9644
$ws.Cells.Item(64, 10).Value = 7599.8
$ws.Cells.Item(64, 11).Value = 6504.9644
$ws.Cells.Item(64, 12).Value = 7599.8
$ws.Cells.Item(64, 13).Value = -6256.9644
$ws.Cells.Item(64, 14).Value = -8095.8
$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 13).Value = $null
$ws.Cells.Item(67, 8).Value = 6670.8486
$ws.Cells.Item(67, 9).Value = 6504.9644
$ws.Cells.Item(67, 10).Value = 7599.8
$ws.Cells.Item(67, 11).Value = 6504.9644
$ws.Cells.Item(67, 12).Value = 7599.8
$ws.Cells.Item(67, 13).Value = -5646.9644
$ws.Cells.Item(67, 14).Value = -9315.799999999999
$ws.Cells.Item(70, 8).Value = 598665.4
$ws.Cells.Item(70, 9).Value = 1011130.2
$ws.Cells.Item(70, 10).Value = 9429.857
$ws.Cells.Item(70, 11).Value = 3033390.6
$ws.Cells.Item(70, 12).Value = 28289.571
$ws.Cells.Item(70, 13).Value = -3033120.6
$ws.Cells.Item(70, 14).Value = -28829.571
$ws.Cells.Item(73, 8).Value = 598665.4
$ws.Cells.Item(73, 9).Value = 1011130.2
$ws.Cells.Item(73, 10).Value = 9429.857
$ws.Cells.Item(73, 11).Value = 3033390.6
$ws.Cells.Item(73, 12).Value = 28289.571
$ws.Cells.Item(73, 13).Value = -3032454.6
$ws.Cells.Item(73, 14).Value = -30161.571
$ws.Cells.Item(92, 8).Value = 492.0357
$ws.Cells.Item(92, 10).Value = 382.625
$ws.Cells.Item(92, 12).Value = 382.625
$ws.Cells.Item(92, 14).Value = -2878.625
$ws.Cells.Item(98, 8).Value = 4931.1665
$ws.Cells.Item(98, 9).Value = 4931.1665
$ws.Cells.Item(98, 11).Value = 4931.1665
$ws.Cells.Item(98, 13).Value = -3433.1665
$ws.Cells.Item(113, 8).Value = 33964404
$ws.Cells.Item(113, 9).Value = 13892159
$ws.Cells.Item(113, 10).Value = 50022196
$ws.Cells.Item(113, 11).Value = 13892159
$ws.Cells.Item(113, 12).Value = 50022196
$ws.Cells.Item(113, 13).Value = -13888905
$ws.Cells.Item(113, 14).Value = -50028704
$ws.Cells.Item(122, 8).Value = 4931.1665
$ws.Cells.Item(122, 9).Value = 4931.1665
$ws.Cells.Item(122, 11).Value = 14793.4995
$ws.Cells.Item(122, 13).Value = -12343.4995
$ws.Cells.Item(132, 8).Value = 1653.6
$ws.Cells.Item(132, 9).Value = 1090.8
$ws.Cells.Item(132, 11).Value = 3272.4
$ws.Cells.Item(132, 13).Value = -742.3999999999996
$ws.Cells.Item(137, 8).Value = 6640.7085
$ws.Cells.Item(137, 9).Value = 4987.875
$ws.Cells.Item(137, 10).Value = 7467.125
$ws.Cells.Item(137, 11).Value = 14963.625
$ws.Cells.Item(137, 12).Value = 22401.375
$ws.Cells.Item(137, 13).Value = -12413.625
$ws.Cells.Item(137, 14).Value = -27501.375
$ws.Cells.Item(138, 8).Value = 1616934
$ws.Cells.Item(138, 9).Value = 2352.1924
$ws.Cells.Item(138, 10).Value = 2783021
$ws.Cells.Item(138, 11).Value = 7056.5772
$ws.Cells.Item(138, 12).Value = 8349063
$ws.Cells.Item(138, 13).Value = -1916.5772
$ws.Cells.Item(138, 14).Value = -8359343
$ws.Cells.Item(141, 8).Value = 5313.737
$ws.Cells.Item(141, 9).Value = 4754.9287
$ws.Cells.Item(141, 10).Value = 6878.4
$ws.Cells.Item(141, 11).Value = 14264.7861
$ws.Cells.Item(141, 12).Value = 20635.2
$ws.Cells.Item(141, 13).Value = -9084.786100000001
$ws.Cells.Item(141, 14).Value = -30995.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 3089.6875
$ws.Cells.Item(2, 9).Value = 2837.0908
$ws.Cells.Item(2, 10).Value = 3645.4
$ws.Cells.Item(2, 11).Value = 2837.0908
$ws.Cells.Item(2, 12).Value = 3645.4
$ws.Cells.Item(2, 13).Value = -2724.0908
$ws.Cells.Item(2, 14).Value = -3871.4
$ws.Cells.Item(32, 8).Value = 2821441.5
$ws.Cells.Item(32, 9).Value = 2944365.5
$ws.Cells.Item(32, 11).Value = 2944365.5
$ws.Cells.Item(32, 13).Value = -2944078.5
$ws.Cells.Item(61, 8).Value = 6991.4165
$ws.Cells.Item(61, 9).Value = 2118.1333
$ws.Cells.Item(61, 11).Value = 2118.1333
$ws.Cells.Item(61, 13).Value = -1906.1333
$ws.Cells.Item(74, 8).Value = 39969.555
$ws.Cells.Item(74, 9).Value = 57454.332
$ws.Cells.Item(74, 10).Value = 5000
$ws.Cells.Item(74, 11).Value = 57454.332
$ws.Cells.Item(74, 12).Value = 5000
$ws.Cells.Item(74, 13).Value = -56580.332
$ws.Cells.Item(74, 14).Value = -6748
$ws.Cells.Item(77, 8).Value = 39969.555
$ws.Cells.Item(77, 9).Value = 57454.332
$ws.Cells.Item(77, 10).Value = 5000
$ws.Cells.Item(77, 11).Value = 287271.66
$ws.Cells.Item(77, 12).Value = 25000
$ws.Cells.Item(77, 13).Value = -282903.66
$ws.Cells.Item(77, 14).Value = -33736
$ws.Cells.Item(97, 8).Value = 25004150
$ws.Cells.Item(97, 9).Value = 27782192
$ws.Cells.Item(97, 11).Value = 27782192
$ws.Cells.Item(97, 13).Value = -27781696
$ws.Cells.Item(115, 8).Value = 59688
$ws.Cells.Item(115, 10).Value = 59688
$ws.Cells.Item(115, 12).Value = 59688
$ws.Cells.Item(115, 14).Value = -62822
$ws.Cells.Item(116, 8).Value = 3089.6875
$ws.Cells.Item(116, 9).Value = 2837.0908
$ws.Cells.Item(116, 10).Value = 3645.4
$ws.Cells.Item(116, 11).Value = 2837.0908
$ws.Cells.Item(116, 12).Value = 3645.4
$ws.Cells.Item(116, 13).Value = -543.0907999999999
$ws.Cells.Item(116, 14).Value = -8233.4
$ws.Cells.Item(122, 8).Value = 4318.2173
$ws.Cells.Item(122, 9).Value = 1954.1111
$ws.Cells.Item(122, 10).Value = 5838
$ws.Cells.Item(122, 11).Value = 5862.3333
$ws.Cells.Item(122, 12).Value = 17514
$ws.Cells.Item(122, 13).Value = -3412.3333
$ws.Cells.Item(122, 14).Value = -22414
$ws.Cells.Item(124, 8).Value = 61041
$ws.Cells.Item(124, 10).Value = 61041
$ws.Cells.Item(124, 12).Value = 61041
$ws.Cells.Item(124, 14).Value = -70861
$ws.Cells.Item(131, 8).Value = 51586.5
$ws.Cells.Item(131, 10).Value = 51586.5
$ws.Cells.Item(131, 12).Value = 51586.5
$ws.Cells.Item(131, 14).Value = -61666.5
$ws.Cells.Item(132, 8).Value = 8887.306
$ws.Cells.Item(132, 9).Value = 9204.5
$ws.Cells.Item(132, 11).Value = 27613.5
$ws.Cells.Item(132, 13).Value = -25083.5
$ws.Cells.Item(136, 8).Value = 6991.4165
$ws.Cells.Item(136, 9).Value = 2118.1333
$ws.Cells.Item(136, 11).Value = 6354.3999
$ws.Cells.Item(136, 13).Value = -3804.3999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 3089.6875
$ws.Cells.Item(3, 9).Value = 2837.0908
$ws.Cells.Item(3, 10).Value = 3645.4
$ws.Cells.Item(3, 11).Value = 2837.0908
$ws.Cells.Item(3, 12).Value = 3645.4
$ws.Cells.Item(3, 13).Value = -2723.0908
$ws.Cells.Item(3, 14).Value = -3873.4
$ws.Cells.Item(5, 8).Value = 1218.2
$ws.Cells.Item(5, 9).Value = 496
$ws.Cells.Item(5, 10).Value = 1699.6666
$ws.Cells.Item(5, 11).Value = 496
$ws.Cells.Item(5, 12).Value = 1699.6666
$ws.Cells.Item(5, 13).Value = -383
$ws.Cells.Item(5, 14).Value = -1925.6666
$ws.Cells.Item(20, 8).Value = 44876930
$ws.Cells.Item(20, 9).Value = 50727876
$ws.Cells.Item(20, 11).Value = 50727876
$ws.Cells.Item(20, 13).Value = -50727629
$ws.Cells.Item(86, 8).Value = 30035.344
$ws.Cells.Item(86, 9).Value = 40699.08
$ws.Cells.Item(86, 10).Value = 3376
$ws.Cells.Item(86, 11).Value = 40699.08
$ws.Cells.Item(86, 12).Value = 3376
$ws.Cells.Item(86, 13).Value = -39576.08
$ws.Cells.Item(86, 14).Value = -5622
$ws.Cells.Item(89, 8).Value = 30035.344
$ws.Cells.Item(89, 9).Value = 40699.08
$ws.Cells.Item(89, 10).Value = 3376
$ws.Cells.Item(89, 11).Value = 203495.4
$ws.Cells.Item(89, 12).Value = 16880
$ws.Cells.Item(89, 13).Value = -197879.4
$ws.Cells.Item(89, 14).Value = -28112
$ws.Cells.Item(94, 8).Value = 3302.9285
$ws.Cells.Item(94, 9).Value = 2221.4546
$ws.Cells.Item(94, 10).Value = 7268.3335
$ws.Cells.Item(94, 11).Value = 2221.4546
$ws.Cells.Item(94, 12).Value = 7268.3335
$ws.Cells.Item(94, 13).Value = -1770.4546
$ws.Cells.Item(94, 14).Value = -8170.3335
$ws.Cells.Item(105, 8).Value = 5869.3
$ws.Cells.Item(105, 9).Value = 6594.25
$ws.Cells.Item(105, 11).Value = 6594.25
$ws.Cells.Item(105, 13).Value = -4847.25
$ws.Cells.Item(107, 8).Value = 112511450
$ws.Cells.Item(107, 9).Value = 125012160
$ws.Cells.Item(107, 11).Value = 125012160
$ws.Cells.Item(107, 13).Value = -125010240
$ws.Cells.Item(134, 8).Value = 5275.196
$ws.Cells.Item(134, 9).Value = 2353.6775
$ws.Cells.Item(134, 10).Value = 11313
$ws.Cells.Item(134, 11).Value = 7061.032499999999
$ws.Cells.Item(134, 12).Value = 33939
$ws.Cells.Item(134, 13).Value = -4526.032499999999
$ws.Cells.Item(134, 14).Value = -39009

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 6177.1035
$ws.Cells.Item(16, 9).Value = 733.75
$ws.Cells.Item(16, 10).Value = 8250.762000000001
$ws.Cells.Item(16, 11).Value = 733.75
$ws.Cells.Item(16, 12).Value = 8250.762000000001
$ws.Cells.Item(16, 13).Value = -446.75
$ws.Cells.Item(16, 14).Value = -8824.762000000001
$ws.Cells.Item(22, 8).Value = 500
$ws.Cells.Item(22, 9).Value = 500
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 500
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -150
$ws.Cells.Item(22, 14).Value = $null
$ws.Cells.Item(31, 8).Value = 7775.95
$ws.Cells.Item(31, 9).Value = 2185
$ws.Cells.Item(31, 11).Value = 2185
$ws.Cells.Item(31, 13).Value = -1890
$ws.Cells.Item(34, 8).Value = 7775.95
$ws.Cells.Item(34, 9).Value = 2185
$ws.Cells.Item(34, 11).Value = 2185
$ws.Cells.Item(34, 13).Value = -1983
$ws.Cells.Item(47, 8).Value = 1000000
$ws.Cells.Item(47, 10).Value = 1000000
$ws.Cells.Item(47, 12).Value = 1000000
$ws.Cells.Item(47, 14).Value = -1001132
$ws.Cells.Item(58, 8).Value = 5282.9414
$ws.Cells.Item(58, 10).Value = 8778.4
$ws.Cells.Item(58, 12).Value = 8778.4
$ws.Cells.Item(58, 14).Value = -9184.4
$ws.Cells.Item(62, 8).Value = 24310220
$ws.Cells.Item(62, 10).Value = 7498.3335
$ws.Cells.Item(62, 12).Value = 7498.3335
$ws.Cells.Item(62, 14).Value = -8746.333500000001
$ws.Cells.Item(65, 8).Value = 24310220
$ws.Cells.Item(65, 10).Value = 7498.3335
$ws.Cells.Item(65, 12).Value = 37491.6675
$ws.Cells.Item(65, 14).Value = -43731.6675
$ws.Cells.Item(70, 8).Value = 40000
$ws.Cells.Item(70, 9).Value = 40000
$ws.Cells.Item(70, 11).Value = 40000
$ws.Cells.Item(70, 13).Value = -39685
$ws.Cells.Item(73, 8).Value = 40000
$ws.Cells.Item(73, 9).Value = 40000
$ws.Cells.Item(73, 11).Value = 40000
$ws.Cells.Item(73, 13).Value = -38908
$ws.Cells.Item(105, 8).Value = 5953150.5
$ws.Cells.Item(105, 9).Value = 6494119
$ws.Cells.Item(105, 11).Value = 6494119
$ws.Cells.Item(105, 13).Value = -6492372
$ws.Cells.Item(113, 8).Value = 6177.1035
$ws.Cells.Item(113, 9).Value = 733.75
$ws.Cells.Item(113, 10).Value = 8250.762000000001
$ws.Cells.Item(113, 11).Value = 733.75
$ws.Cells.Item(113, 12).Value = 8250.762000000001
$ws.Cells.Item(113, 13).Value = 1436.25
$ws.Cells.Item(113, 14).Value = -12590.762
$ws.Cells.Item(119, 8).Value = 124880.5
$ws.Cells.Item(119, 10).Value = 124880.5
$ws.Cells.Item(119, 12).Value = 124880.5
$ws.Cells.Item(119, 14).Value = -134556.5
$ws.Cells.Item(122, 8).Value = 2670.2727
$ws.Cells.Item(122, 9).Value = 1789.8572
$ws.Cells.Item(122, 11).Value = 5369.571599999999
$ws.Cells.Item(122, 13).Value = -2919.571599999999
$ws.Cells.Item(134, 8).Value = 5576.081
$ws.Cells.Item(134, 9).Value = 1770.0625
$ws.Cells.Item(134, 11).Value = 5310.1875
$ws.Cells.Item(134, 13).Value = -2775.1875
$ws.Cells.Item(136, 8).Value = 5282.9414
$ws.Cells.Item(136, 10).Value = 8778.4
$ws.Cells.Item(136, 12).Value = 26335.2
$ws.Cells.Item(136, 14).Value = -31435.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 42723370
$ws.Cells.Item(4, 9).Value = 78652080
$ws.Cells.Item(4, 10).Value = 806542.4399999999
$ws.Cells.Item(4, 11).Value = 235956240
$ws.Cells.Item(4, 12).Value = 2419627.32
$ws.Cells.Item(4, 13).Value = -235956128
$ws.Cells.Item(4, 14).Value = -2419851.32
$ws.Cells.Item(17, 8).Value = 1236.1765
$ws.Cells.Item(17, 10).Value = 2513.5715
$ws.Cells.Item(17, 12).Value = 7540.7145
$ws.Cells.Item(17, 14).Value = -7878.7145
$ws.Cells.Item(34, 8).Value = 4915
$ws.Cells.Item(34, 10).Value = 5375.7617
$ws.Cells.Item(34, 12).Value = 16127.2851
$ws.Cells.Item(34, 14).Value = -16295.2851
$ws.Cells.Item(39, 8).Value = 8374.923000000001
$ws.Cells.Item(39, 9).Value = 5799.2
$ws.Cells.Item(39, 10).Value = 9984.75
$ws.Cells.Item(39, 11).Value = 17397.6
$ws.Cells.Item(39, 12).Value = 29954.25
$ws.Cells.Item(39, 13).Value = -17103.6
$ws.Cells.Item(39, 14).Value = -30542.25
$ws.Cells.Item(62, 8).Value = 2875
$ws.Cells.Item(62, 9).Value = 2875
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 11).Value = 8625
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 13).Value = -7939
$ws.Cells.Item(62, 14).Value = $null
$ws.Cells.Item(65, 8).Value = 2875
$ws.Cells.Item(65, 9).Value = 2875
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 25875
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).Value = -22443
$ws.Cells.Item(65, 14).Value = $null
$ws.Cells.Item(120, 8).Value = 16308.546
$ws.Cells.Item(120, 9).Value = 5848.75
$ws.Cells.Item(120, 10).Value = 22285.572
$ws.Cells.Item(120, 11).Value = 17546.25
$ws.Cells.Item(120, 12).Value = 66856.716
$ws.Cells.Item(120, 13).Value = -12708.25
$ws.Cells.Item(120, 14).Value = -76532.716
$ws.Cells.Item(123, 8).Value = 2300
$ws.Cells.Item(123, 9).Value = 2300
$ws.Cells.Item(123, 10).Value = 0
$ws.Cells.Item(123, 11).Value = 6900
$ws.Cells.Item(123, 12).Value = 0
$ws.Cells.Item(123, 13).Value = -4450
$ws.Cells.Item(123, 14).Value = $null
$ws.Cells.Item(132, 8).Value = 4724.2793
$ws.Cells.Item(132, 9).Value = 2253.4666
$ws.Cells.Item(132, 10).Value = 6047.9287
$ws.Cells.Item(132, 11).Value = 20281.1994
$ws.Cells.Item(132, 12).Value = 54431.35830000001
$ws.Cells.Item(132, 13).Value = -17751.1994
$ws.Cells.Item(132, 14).Value = -59491.35830000001
$ws.Cells.Item(139, 8).Value = 57583.473
$ws.Cells.Item(139, 9).Value = 79545.16
$ws.Cells.Item(139, 11).Value = 238635.48
$ws.Cells.Item(139, 13).Value = -233495.48
$ws.Cells.Item(141, 8).Value = 4772.2
$ws.Cells.Item(141, 9).Value = 3512.0908
$ws.Cells.Item(141, 10).Value = 8237.5
$ws.Cells.Item(141, 11).Value = 10536.2724
$ws.Cells.Item(141, 12).Value = 24712.5
$ws.Cells.Item(141, 13).Value = -5356.2724
$ws.Cells.Item(141, 14).Value = -35072.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(26, 8).Value = 24950
$ws.Cells.Item(26, 10).Value = 24950
$ws.Cells.Item(26, 12).Value = 24950
$ws.Cells.Item(26, 14).Value = -25510
$ws.Cells.Item(49, 8).Value = 13234.25
$ws.Cells.Item(50, 8).Value = 24950
$ws.Cells.Item(50, 10).Value = 24950
$ws.Cells.Item(50, 12).Value = 24950
$ws.Cells.Item(50, 14).Value = -25946
$ws.Cells.Item(80, 8).Value = 2342
$ws.Cells.Item(80, 9).Value = 2342
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 11).Value = 2342
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 13).Value = -1344
$ws.Cells.Item(80, 14).Value = $null
$ws.Cells.Item(83, 8).Value = 2342
$ws.Cells.Item(83, 9).Value = 2342
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 11).Value = 11710
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 13).Value = -6718
$ws.Cells.Item(83, 14).Value = $null
$ws.Cells.Item(102, 8).Value = 4358.8887
$ws.Cells.Item(102, 9).Value = 4129.087
$ws.Cells.Item(102, 10).Value = 5680.25
$ws.Cells.Item(102, 11).Value = 4129.087
$ws.Cells.Item(102, 12).Value = 5680.25
$ws.Cells.Item(102, 13).Value = -2507.087
$ws.Cells.Item(102, 14).Value = -8924.25
$ws.Cells.Item(113, 8).Value = 7626.645
$ws.Cells.Item(113, 9).Value = 4785.643
$ws.Cells.Item(113, 11).Value = 4785.643
$ws.Cells.Item(113, 13).Value = -2615.643
$ws.Cells.Item(122, 8).Value = 77551.19
$ws.Cells.Item(122, 9).Value = 136153.44
$ws.Cells.Item(122, 10).Value = 2205.4285
$ws.Cells.Item(122, 11).Value = 408460.32
$ws.Cells.Item(122, 12).Value = 6616.2855
$ws.Cells.Item(122, 13).Value = -406010.32
$ws.Cells.Item(122, 14).Value = -11516.2855
$ws.Cells.Item(123, 8).Value = 30000
$ws.Cells.Item(123, 10).Value = 30000
$ws.Cells.Item(123, 12).Value = 30000
$ws.Cells.Item(123, 14).Value = -34900
$ws.Cells.Item(126, 8).Value = 5091.857
$ws.Cells.Item(126, 9).Value = 2642.5715
$ws.Cells.Item(126, 10).Value = 7541.143
$ws.Cells.Item(126, 11).Value = 7927.7145
$ws.Cells.Item(126, 12).Value = 22623.429
$ws.Cells.Item(126, 13).Value = -5457.7145
$ws.Cells.Item(126, 14).Value = -27563.429
$ws.Cells.Item(132, 8).Value = 4903.385
$ws.Cells.Item(132, 9).Value = 2642.8096
$ws.Cells.Item(132, 10).Value = 14397.8
$ws.Cells.Item(132, 11).Value = 7928.4288
$ws.Cells.Item(132, 12).Value = 43193.39999999999
$ws.Cells.Item(132, 13).Value = -5398.4288
$ws.Cells.Item(132, 14).Value = -48253.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6503.3335
$ws.Cells.Item(7, 9).Value = 5036.2
$ws.Cells.Item(7, 10).Value = 8337.25
$ws.Cells.Item(7, 11).Value = 5036.2
$ws.Cells.Item(7, 12).Value = 8337.25
$ws.Cells.Item(7, 13).Value = -4924.2
$ws.Cells.Item(7, 14).Value = -8561.25
$ws.Cells.Item(40, 8).Value = 4537.091
$ws.Cells.Item(40, 9).Value = 2541.6
$ws.Cells.Item(40, 10).Value = 6200
$ws.Cells.Item(40, 11).Value = 2541.6
$ws.Cells.Item(40, 12).Value = 6200
$ws.Cells.Item(40, 13).Value = -2405.6
$ws.Cells.Item(40, 14).Value = -6472
$ws.Cells.Item(42, 8).Value = 16500
$ws.Cells.Item(42, 9).Value = 8000
$ws.Cells.Item(42, 11).Value = 8000
$ws.Cells.Item(42, 13).Value = -7437
$ws.Cells.Item(49, 8).Value = 16500
$ws.Cells.Item(49, 9).Value = 8000
$ws.Cells.Item(49, 11).Value = 8000
$ws.Cells.Item(49, 13).Value = -7853
$ws.Cells.Item(68, 8).Value = 6311.2144
$ws.Cells.Item(68, 9).Value = 4489.5
$ws.Cells.Item(68, 10).Value = 7039.9
$ws.Cells.Item(68, 11).Value = 4489.5
$ws.Cells.Item(68, 12).Value = 7039.9
$ws.Cells.Item(68, 13).Value = -3740.5
$ws.Cells.Item(68, 14).Value = -8537.9
$ws.Cells.Item(71, 8).Value = 6311.2144
$ws.Cells.Item(71, 9).Value = 4489.5
$ws.Cells.Item(71, 10).Value = 7039.9
$ws.Cells.Item(71, 11).Value = 22447.5
$ws.Cells.Item(71, 12).Value = 35199.5
$ws.Cells.Item(71, 13).Value = -18703.5
$ws.Cells.Item(71, 14).Value = -42687.5
$ws.Cells.Item(82, 8).Value = 59413176
$ws.Cells.Item(82, 10).Value = 2495.7144
$ws.Cells.Item(82, 12).Value = 2495.7144
$ws.Cells.Item(82, 14).Value = -3217.7144
$ws.Cells.Item(85, 8).Value = 59413176
$ws.Cells.Item(85, 10).Value = 2495.7144
$ws.Cells.Item(85, 12).Value = 2495.7144
$ws.Cells.Item(85, 14).Value = -4991.7144
$ws.Cells.Item(122, 8).Value = 4387.4287
$ws.Cells.Item(122, 9).Value = 2704
$ws.Cells.Item(122, 11).Value = 8112
$ws.Cells.Item(122, 13).Value = -5662
$ws.Cells.Item(126, 8).Value = 6503.3335
$ws.Cells.Item(126, 9).Value = 5036.2
$ws.Cells.Item(126, 10).Value = 8337.25
$ws.Cells.Item(126, 11).Value = 15108.6
$ws.Cells.Item(126, 12).Value = 25011.75
$ws.Cells.Item(126, 13).Value = -12638.6
$ws.Cells.Item(126, 14).Value = -29951.75
$ws.Cells.Item(127, 8).Value = 47384.11
$ws.Cells.Item(127, 10).Value = 47384.11
$ws.Cells.Item(127, 12).Value = 47384.11
$ws.Cells.Item(127, 14).Value = -57304.11
$ws.Cells.Item(132, 8).Value = 5579.82
$ws.Cells.Item(132, 9).Value = 2721.889
$ws.Cells.Item(132, 10).Value = 8934.781999999999
$ws.Cells.Item(132, 11).Value = 8165.667
$ws.Cells.Item(132, 12).Value = 26804.346
$ws.Cells.Item(132, 13).Value = -5635.667
$ws.Cells.Item(132, 14).Value = -31864.346
$ws.Cells.Item(136, 8).Value = 13016.523
$ws.Cells.Item(136, 9).Value = 3221
$ws.Cells.Item(136, 10).Value = 21108.479
$ws.Cells.Item(136, 11).Value = 9663
$ws.Cells.Item(136, 12).Value = 63325.437
$ws.Cells.Item(136, 13).Value = -7113
$ws.Cells.Item(136, 14).Value = -68425.43700000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 1011.6667
$ws.Cells.Item(100, 9).Value = 332.5
$ws.Cells.Item(100, 10).Value = 1205.7142
$ws.Cells.Item(100, 11).Value = 665
$ws.Cells.Item(100, 12).Value = 2411.4284
$ws.Cells.Item(100, 13).Value = -124
$ws.Cells.Item(100, 14).Value = -3493.4284
$ws.Cells.Item(113, 8).Value = 1033.8823
$ws.Cells.Item(113, 9).Value = 1182.5625
$ws.Cells.Item(113, 11).Value = 3547.6875
$ws.Cells.Item(113, 13).Value = -1377.6875
$ws.Cells.Item(122, 8).Value = 3581.7932
$ws.Cells.Item(122, 9).Value = 2753.9546
$ws.Cells.Item(122, 10).Value = 6183.5713
$ws.Cells.Item(122, 11).Value = 8261.863799999999
$ws.Cells.Item(122, 12).Value = 18550.7139
$ws.Cells.Item(122, 13).Value = -5811.863799999999
$ws.Cells.Item(122, 14).Value = -23450.7139
$ws.Cells.Item(132, 8).Value = 8563
$ws.Cells.Item(132, 9).Value = 17732.5
$ws.Cells.Item(132, 11).Value = 53197.5
$ws.Cells.Item(132, 13).Value = -50667.5
$ws.Cells.Item(136, 8).Value = 4542.4165
$ws.Cells.Item(136, 9).Value = 1678.3334
$ws.Cells.Item(136, 11).Value = 5035.0002
$ws.Cells.Item(136, 13).Value = -2485.0002
